$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the numeric-looking text columns (Donation Amount, Phone, Zip) as Text
# *before* writing values, so Excel stores them as strings (matching the source
# data format used by cells B2:L3) instead of auto-converting to numbers.
$ws.Range("F4:F12").NumberFormat = "@"
$ws.Range("H4:H12").NumberFormat = "@"
$ws.Range("L4:L12").NumberFormat = "@"

# Row 4
$ws.Range("A4").Value = 'REC-1741876652886-252'
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = '2025-03-13T14:37:32.907Z'
$ws.Range("D4").Value = 'Deepak'
$ws.Range("E4").Value = 'Adhikari'
$ws.Range("F4").Value = '3445'
$ws.Range("G4").Value = 'dadhikari856@gmail.com'
$ws.Range("H4").Value = '3477712375'
$ws.Range("I4").Value = '11 alpine ln'
$ws.Range("J4").Value = 'Hicksville'
$ws.Range("K4").Value = 'NY'
$ws.Range("L4").Value = '11801'

# Row 5
$ws.Range("A5").Value = 'REC-1741876747020-229'
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = '2025-03-13T14:39:07.031Z'
$ws.Range("D5").Value = 'Deepak'
$ws.Range("E5").Value = 'Adhikari'
$ws.Range("F5").Value = '3445'
$ws.Range("G5").Value = 'dadhikari856@gmail.com'
$ws.Range("H5").Value = '3477712375'
$ws.Range("I5").Value = '11 alpine ln'
$ws.Range("J5").Value = 'Hicksville'
$ws.Range("K5").Value = 'NY'
$ws.Range("L5").Value = '11801'

# Row 6
$ws.Range("A6").Value = 'REC-1741876748126-582'
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = '2025-03-13T14:39:08.134Z'
$ws.Range("D6").Value = 'Deepak'
$ws.Range("E6").Value = 'Adhikari'
$ws.Range("F6").Value = '3445'
$ws.Range("G6").Value = 'dadhikari856@gmail.com'
$ws.Range("H6").Value = '3477712375'
$ws.Range("I6").Value = '11 alpine ln'
$ws.Range("J6").Value = 'Hicksville'
$ws.Range("K6").Value = 'NY'
$ws.Range("L6").Value = '11801'

# Row 7
$ws.Range("A7").Value = 'REC-1741876748336-135'
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = '2025-03-13T14:39:08.338Z'
$ws.Range("D7").Value = 'Deepak'
$ws.Range("E7").Value = 'Adhikari'
$ws.Range("F7").Value = '3445'
$ws.Range("G7").Value = 'dadhikari856@gmail.com'
$ws.Range("H7").Value = '3477712375'
$ws.Range("I7").Value = '11 alpine ln'
$ws.Range("J7").Value = 'Hicksville'
$ws.Range("K7").Value = 'NY'
$ws.Range("L7").Value = '11801'

# Row 8
$ws.Range("A8").Value = 'REC-1741876873473-419'
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = '2025-03-13T14:41:13.475Z'
$ws.Range("D8").Value = 'Deepak'
$ws.Range("E8").Value = 'Adhikari'
$ws.Range("F8").Value = '3445'
$ws.Range("G8").Value = 'dadhikari856@gmail.com'
$ws.Range("H8").Value = '3477712375'
$ws.Range("I8").Value = '11 alpine ln'
$ws.Range("J8").Value = 'Hicksville'
$ws.Range("K8").Value = 'NY'
$ws.Range("L8").Value = '11801'

# Row 9
$ws.Range("A9").Value = 'REC-1741876962971-397'
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = '2025-03-13T14:42:42.973Z'
$ws.Range("D9").Value = 'Deepak'
$ws.Range("E9").Value = 'Adhikari'
$ws.Range("F9").Value = '3445'
$ws.Range("G9").Value = 'dadhikari856@gmail.com'
$ws.Range("H9").Value = '3477712375'
$ws.Range("I9").Value = '11 alpine ln'
$ws.Range("J9").Value = 'Hicksville'
$ws.Range("K9").Value = 'NY'
$ws.Range("L9").Value = '11801'

# Row 10
$ws.Range("A10").Value = 'REC-1741876992736-933'
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = '2025-03-13T14:43:12.740Z'
$ws.Range("D10").Value = 'Deepak'
$ws.Range("E10").Value = 'Adhikari'
$ws.Range("F10").Value = '3445'
$ws.Range("G10").Value = 'dadhikari856@gmail.com'
$ws.Range("H10").Value = '3477712375'
$ws.Range("I10").Value = '11 alpine ln'
$ws.Range("J10").Value = 'Hicksville'
$ws.Range("K10").Value = 'NY'
$ws.Range("L10").Value = '11801'

# Row 11
$ws.Range("A11").Value = 'REC-1741877009123-108'
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = '2025-03-13T14:43:29.132Z'
$ws.Range("D11").Value = 'Deepak'
$ws.Range("E11").Value = 'Adhikari'
$ws.Range("F11").Value = '3445'
$ws.Range("G11").Value = 'dadhikari856@gmail.com'
$ws.Range("H11").Value = '3477712375'
$ws.Range("I11").Value = '11 alpine ln'
$ws.Range("J11").Value = 'Hicksville'
$ws.Range("K11").Value = 'NY'
$ws.Range("L11").Value = '11801'

# Row 12
$ws.Range("A12").Value = 'REC-1741877081275-061'
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = '2025-03-13T14:44:41.278Z'
$ws.Range("D12").Value = 'Deepak'
$ws.Range("E12").Value = 'Adhikari'
$ws.Range("F12").Value = '3445'
$ws.Range("G12").Value = 'dadhikari856@gmail.com'
$ws.Range("H12").Value = '3477712375'
$ws.Range("I12").Value = '11 alpine ln'
$ws.Range("J12").Value = 'Hicksville'
$ws.Range("K12").Value = 'NY'
$ws.Range("L12").Value = '11801'
